$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Pre-format the two new rows by copying row 2's look (borders/fill/font/number format) ---
$ws.Range("A2:G2").Copy() | Out-Null
$ws.Range("A3:G3").PasteSpecial(-4122) | Out-Null
$ws.Range("A2:G2").Copy() | Out-Null
$ws.Range("A4:G4").PasteSpecial(-4122) | Out-Null

# --- New row 3: ORG_PV_Phone scenario ---
$ws.Range("A3").Value = "0000_ORG_PV_Phone_Scenario"
$ws.Range("B3").Value = "Verify the ORG_PV_Phone Info"
$ws.Range("C3").Value = "Verify_All_Buckets_ORG_PV_PHONE"

# --- Update existing row 2: the R3 row-count / Execution flag for the Phone scenario ---
$ws.Range("E2").Value = "No"
$ws.Range("D2").Value = "5"

# --- New row 4: ORG_Website_Cache scenario ---
$ws.Range("A4").Value = "0000_ORG_Website_Cache_Scenario"
$ws.Range("B4").Value = "Verify the ORG_Website_Cache_ Info"
$ws.Range("C4").Value = "Verify_ORG_WEBSITE_CACHE"
$ws.Range("D4").Value = "1"

# --- Finish row 3 ---
$ws.Range("D3").Value = "100"
$ws.Range("E3").Value = "Yes"

# --- Finish row 4 ---
$ws.Range("E4").Value = "No"

# --- Remaining static columns for both new rows ---
$ws.Range("F3").Value = "End-To-End_Testcases"
$ws.Range("G3").Value = "Selvamani_M"
$ws.Range("F4").Value = "End-To-End_Testcases"
$ws.Range("G4").Value = "Selvamani_M"

# --- Widen column A to fit the new, longer scenario names ---
$ws.Columns.Item(1).ColumnWidth = 39.3

# --- Move the active selection, matching the author's last cursor position ---
$ws.Range("D19").Select() | Out-Null
